$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Adam12"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = [double]"4.480363666666666"
$ws.Cells.Item(2, 8).Value = [double]"13.441091"
$ws.Cells.Item(2, 9).Value = [double]"0.05823429740900917"
$ws.Cells.Item(2, 10).Value = [double]"0.05886574272937452"
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = [double]"121.928739"
$ws.Cells.Item(2, 14).Value = [double]"365.786217"
$ws.Cells.Item(2, 15).Value = [double]"0.2282232151508951"
$ws.Cells.Item(2, 16).Value = [double]"0.2419720431319445"
$ws.Cells.Item(2, 17).Value = [double]"546.285092138083"
$ws.Cells.Item(2, 18).Value = [double]"4916.565829242747"
$ws.Cells.Item(2, 19).Value = [double]"0.01329041858673751"
$ws.Cells.Item(2, 20).Value = [double]"0.01424386403870616"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Adam12"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = [double]"4.480363666666666"
$ws.Cells.Item(3, 8).Value = [double]"13.441091"
$ws.Cells.Item(3, 9).Value = [double]"0.05823429740900917"
$ws.Cells.Item(3, 10).Value = [double]"0.05886574272937452"
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = [double]"147.91433"
$ws.Cells.Item(3, 14).Value = [double]"443.74299"
$ws.Cells.Item(3, 15).Value = [double]"0.2768624053389947"
$ws.Cells.Item(3, 16).Value = [double]"0.2935413991166814"
$ws.Cells.Item(3, 17).Value = [double]"662.7099899113432"
$ws.Cells.Item(3, 18).Value = [double]"5964.389909202089"
$ws.Cells.Item(3, 19).Value = [double]"0.01612288765388467"
$ws.Cells.Item(3, 20).Value = [double]"0.01727953248082321"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Adam12"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = [double]"4.480363666666666"
$ws.Cells.Item(4, 8).Value = [double]"13.441091"
$ws.Cells.Item(4, 9).Value = [double]"0.05823429740900917"
$ws.Cells.Item(4, 10).Value = [double]"0.05886574272937452"
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = [double]"83.50496933333334"
$ws.Cells.Item(4, 14).Value = [double]"250.514908"
$ws.Cells.Item(4, 15).Value = [double]"0.1563025480180701"
$ws.Cells.Item(4, 16).Value = [double]"0.1657186665504434"
$ws.Cells.Item(4, 17).Value = [double]"374.1326305871809"
$ws.Cells.Item(4, 18).Value = [double]"3367.193675284628"
$ws.Cells.Item(4, 19).Value = [double]"0.009102169067070228"
$ws.Cells.Item(4, 20).Value = [double]"0.009755152390613404"

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Adam12"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = [double]"4.480363666666666"
$ws.Cells.Item(5, 8).Value = [double]"13.441091"
$ws.Cells.Item(5, 9).Value = [double]"0.05823429740900917"
$ws.Cells.Item(5, 10).Value = [double]"0.05886574272937452"
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = [double]"91.06846250000001"
$ws.Cells.Item(5, 14).Value = [double]"182.136925"
$ws.Cells.Item(5, 15).Value = [double]"0.1704597085236707"
$ws.Cells.Item(5, 16).Value = [double]"0.1204857969594293"
$ws.Cells.Item(5, 17).Value = [double]"408.0198305641958"
$ws.Cells.Item(5, 18).Value = [double]"2448.118983385175"
$ws.Cells.Item(5, 19).Value = [double]"0.009926601362420457"
$ws.Cells.Item(5, 20).Value = [double]"0.007092485926357422"

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Adam12"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = [double]"4.480363666666666"
$ws.Cells.Item(6, 8).Value = [double]"13.441091"
$ws.Cells.Item(6, 9).Value = [double]"0.05823429740900917"
$ws.Cells.Item(6, 10).Value = [double]"0.05886574272937452"
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = [double]"89.83563"
$ws.Cells.Item(6, 14).Value = [double]"269.50689"
$ws.Cells.Item(6, 15).Value = [double]"0.1681521229683693"
$ws.Cells.Item(6, 16).Value = [double]"0.1782820942415013"
$ws.Cells.Item(6, 17).Value = [double]"402.49629262411"
$ws.Cells.Item(6, 18).Value = [double]"3622.46663361699"
$ws.Cells.Item(6, 19).Value = [double]"0.009792220738896299"
$ws.Cells.Item(6, 20).Value = [double]"0.01049470789287432"

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Adam12"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = [double]"69.95512000000001"
$ws.Cells.Item(7, 8).Value = [double]"209.86536"
$ws.Cells.Item(7, 9).Value = [double]"0.90925370493279"
$ws.Cells.Item(7, 10).Value = [double]"0.9191129120074827"
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = [double]"121.928739"
$ws.Cells.Item(7, 14).Value = [double]"365.786217"
$ws.Cells.Item(7, 15).Value = [double]"0.2282232151508951"
$ws.Cells.Item(7, 16).Value = [double]"0.2419720431319445"
$ws.Cells.Item(7, 17).Value = [double]"8529.53956819368"
$ws.Cells.Item(7, 18).Value = [double]"76765.85611374311"
$ws.Cells.Item(7, 19).Value = [double]"0.2075128039276246"
$ws.Cells.Item(7, 20).Value = [double]"0.2223996291874017"

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Adam12"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = [double]"69.95512000000001"
$ws.Cells.Item(8, 8).Value = [double]"209.86536"
$ws.Cells.Item(8, 9).Value = [double]"0.90925370493279"
$ws.Cells.Item(8, 10).Value = [double]"0.9191129120074827"
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = [double]"147.91433"
$ws.Cells.Item(8, 14).Value = [double]"443.74299"
$ws.Cells.Item(8, 15).Value = [double]"0.2768624053389947"
$ws.Cells.Item(8, 16).Value = [double]"0.2935413991166814"
$ws.Cells.Item(8, 17).Value = [double]"10347.3647048696"
$ws.Cells.Item(8, 18).Value = [double]"93126.2823438264"
$ws.Cells.Item(8, 19).Value = [double]"0.2517381678110848"
$ws.Cells.Item(8, 20).Value = [double]"0.2697976901368838"

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Adam12"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = [double]"69.95512000000001"
$ws.Cells.Item(9, 8).Value = [double]"209.86536"
$ws.Cells.Item(9, 9).Value = [double]"0.90925370493279"
$ws.Cells.Item(9, 10).Value = [double]"0.9191129120074827"
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = [double]"83.50496933333334"
$ws.Cells.Item(9, 14).Value = [double]"250.514908"
$ws.Cells.Item(9, 15).Value = [double]"0.1563025480180701"
$ws.Cells.Item(9, 16).Value = [double]"0.1657186665504434"
$ws.Cells.Item(9, 17).Value = [double]"5841.600150309654"
$ws.Cells.Item(9, 18).Value = [double]"52574.40135278688"
$ws.Cells.Item(9, 19).Value = [double]"0.1421186708758655"
$ws.Cells.Item(9, 20).Value = [double]"0.152314166187175"

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Adam12"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = [double]"69.95512000000001"
$ws.Cells.Item(10, 8).Value = [double]"209.86536"
$ws.Cells.Item(10, 9).Value = [double]"0.90925370493279"
$ws.Cells.Item(10, 10).Value = [double]"0.9191129120074827"
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = [double]"91.06846250000001"
$ws.Cells.Item(10, 14).Value = [double]"182.136925"
$ws.Cells.Item(10, 15).Value = [double]"0.1704597085236707"
$ws.Cells.Item(10, 16).Value = [double]"0.1204857969594293"
$ws.Cells.Item(10, 17).Value = [double]"6370.705222403001"
$ws.Cells.Item(10, 18).Value = [double]"38224.23133441801"
$ws.Cells.Item(10, 19).Value = [double]"0.1549911215169111"
$ws.Cells.Item(10, 20).Value = [double]"0.1107400516989234"

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Adam12"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = [double]"69.95512000000001"
$ws.Cells.Item(11, 8).Value = [double]"209.86536"
$ws.Cells.Item(11, 9).Value = [double]"0.90925370493279"
$ws.Cells.Item(11, 10).Value = [double]"0.9191129120074827"
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = [double]"89.83563"
$ws.Cells.Item(11, 14).Value = [double]"269.50689"
$ws.Cells.Item(11, 15).Value = [double]"0.1681521229683693"
$ws.Cells.Item(11, 16).Value = [double]"0.1782820942415013"
$ws.Cells.Item(11, 17).Value = [double]"6284.4622769256"
$ws.Cells.Item(11, 18).Value = [double]"56560.1604923304"
$ws.Cells.Item(11, 19).Value = [double]"0.1528929408013039"
$ws.Cells.Item(11, 20).Value = [double]"0.1638613747970987"

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Adam12"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12, 7).Value = [double]"0.02020466666666667"
$ws.Cells.Item(12, 8).Value = [double]"0.060614"
$ws.Cells.Item(12, 9).Value = [double]"0.0002626136303332581"
$ws.Cells.Item(12, 10).Value = [double]"0.0002654611987820265"
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = [double]"121.928739"
$ws.Cells.Item(12, 14).Value = [double]"365.786217"
$ws.Cells.Item(12, 15).Value = [double]"0.2282232151508951"
$ws.Cells.Item(12, 16).Value = [double]"0.2419720431319445"
$ws.Cells.Item(12, 17).Value = [double]"2.463529528582"
$ws.Cells.Item(12, 18).Value = [double]"22.171765757238"
$ws.Cells.Item(12, 19).Value = [double]"5.993452705710478e-05"
$ws.Cells.Item(12, 20).Value = [double]"6.423418864154223e-05"

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Adam12"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 7).Value = [double]"0.02020466666666667"
$ws.Cells.Item(13, 8).Value = [double]"0.060614"
$ws.Cells.Item(13, 9).Value = [double]"0.0002626136303332581"
$ws.Cells.Item(13, 10).Value = [double]"0.0002654611987820265"
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = [double]"147.91433"
$ws.Cells.Item(13, 14).Value = [double]"443.74299"
$ws.Cells.Item(13, 15).Value = [double]"0.2768624053389947"
$ws.Cells.Item(13, 16).Value = [double]"0.2935413991166814"
$ws.Cells.Item(13, 17).Value = [double]"2.988559732873333"
$ws.Cells.Item(13, 18).Value = [double]"26.89703759586"
$ws.Cells.Item(13, 19).Value = [double]"7.27078413688714e-05"
$ws.Cells.Item(13, 20).Value = [double]"7.792385170166755e-05"

# Row 14
$ws.Cells.Item(14, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 2).Value = "Adam12"
$ws.Cells.Item(14, 3).Value = "Itgb1"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(14, 7).Value = [double]"0.02020466666666667"
$ws.Cells.Item(14, 8).Value = [double]"0.060614"
$ws.Cells.Item(14, 9).Value = [double]"0.0002626136303332581"
$ws.Cells.Item(14, 10).Value = [double]"0.0002654611987820265"
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = [double]"83.50496933333334"
$ws.Cells.Item(14, 14).Value = [double]"250.514908"
$ws.Cells.Item(14, 15).Value = [double]"0.1563025480180701"
$ws.Cells.Item(14, 16).Value = [double]"0.1657186665504434"
$ws.Cells.Item(14, 17).Value = [double]"1.687190070390222"
$ws.Cells.Item(14, 18).Value = [double]"15.184710633512"
$ws.Cells.Item(14, 19).Value = [double]"4.104717956536377e-05"
$ws.Cells.Item(14, 20).Value = [double]"4.399187588303962e-05"

# Row 15
$ws.Cells.Item(15, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 2).Value = "Adam12"
$ws.Cells.Item(15, 3).Value = "Itgb1"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(15, 7).Value = [double]"0.02020466666666667"
$ws.Cells.Item(15, 8).Value = [double]"0.060614"
$ws.Cells.Item(15, 9).Value = [double]"0.0002626136303332581"
$ws.Cells.Item(15, 10).Value = [double]"0.0002654611987820265"
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = [double]"91.06846250000001"
$ws.Cells.Item(15, 14).Value = [double]"182.136925"
$ws.Cells.Item(15, 15).Value = [double]"0.1704597085236707"
$ws.Cells.Item(15, 16).Value = [double]"0.1204857969594293"
$ws.Cells.Item(15, 17).Value = [double]"1.840007928658334"
$ws.Cells.Item(15, 18).Value = [double]"11.04004757195"
$ws.Cells.Item(15, 19).Value = [double]"4.476504288095018e-05"
$ws.Cells.Item(15, 20).Value = [double]"3.198430409705796e-05"

# Row 16
$ws.Cells.Item(16, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 2).Value = "Adam12"
$ws.Cells.Item(16, 3).Value = "Itgb1"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16, 7).Value = [double]"0.02020466666666667"
$ws.Cells.Item(16, 8).Value = [double]"0.060614"
$ws.Cells.Item(16, 9).Value = [double]"0.0002626136303332581"
$ws.Cells.Item(16, 10).Value = [double]"0.0002654611987820265"
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = [double]"89.83563"
$ws.Cells.Item(16, 14).Value = [double]"269.50689"
$ws.Cells.Item(16, 15).Value = [double]"0.1681521229683693"
$ws.Cells.Item(16, 16).Value = [double]"0.1782820942415013"
$ws.Cells.Item(16, 17).Value = [double]"1.81509895894"
$ws.Cells.Item(16, 18).Value = [double]"16.33589063046"
$ws.Cells.Item(16, 19).Value = [double]"4.415903946096788e-05"
$ws.Cells.Item(16, 20).Value = [double]"4.732697845871916e-05"

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Adam12"
$ws.Cells.Item(17, 3).Value = "Itgb1"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = [double]"2.4758755"
$ws.Cells.Item(17, 8).Value = [double]"4.951751"
$ws.Cells.Item(17, 9).Value = [double]"0.03218061767783864"
$ws.Cells.Item(17, 10).Value = [double]"0.02168637206800571"
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = [double]"121.928739"
$ws.Cells.Item(17, 14).Value = [double]"365.786217"
$ws.Cells.Item(17, 15).Value = [double]"0.2282232151508951"
$ws.Cells.Item(17, 16).Value = [double]"0.2419720431319445"
$ws.Cells.Item(17, 17).Value = [double]"301.8803776359945"
$ws.Cells.Item(17, 18).Value = [double]"1811.282265815967"
$ws.Cells.Item(17, 19).Value = [double]"0.007344364031978066"
$ws.Cells.Item(17, 20).Value = [double]"0.005247495757414876"

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Adam12"
$ws.Cells.Item(18, 3).Value = "Itgb1"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = [double]"2.4758755"
$ws.Cells.Item(18, 8).Value = [double]"4.951751"
$ws.Cells.Item(18, 9).Value = [double]"0.03218061767783864"
$ws.Cells.Item(18, 10).Value = [double]"0.02168637206800571"
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = [double]"147.91433"
$ws.Cells.Item(18, 14).Value = [double]"443.74299"
$ws.Cells.Item(18, 15).Value = [double]"0.2768624053389947"
$ws.Cells.Item(18, 16).Value = [double]"0.2935413991166814"
$ws.Cells.Item(18, 17).Value = [double]"366.2174657459149"
$ws.Cells.Item(18, 18).Value = [double]"2197.30479447549"
$ws.Cells.Item(18, 19).Value = [double]"0.00890960321558098"
$ws.Cells.Item(18, 20).Value = [double]"0.006365847998607317"

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Adam12"
$ws.Cells.Item(19, 3).Value = "Itgb1"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = [double]"2.4758755"
$ws.Cells.Item(19, 8).Value = [double]"4.951751"
$ws.Cells.Item(19, 9).Value = [double]"0.03218061767783864"
$ws.Cells.Item(19, 10).Value = [double]"0.02168637206800571"
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = [double]"83.50496933333334"
$ws.Cells.Item(19, 14).Value = [double]"250.514908"
$ws.Cells.Item(19, 15).Value = [double]"0.1563025480180701"
$ws.Cells.Item(19, 16).Value = [double]"0.1657186665504434"
$ws.Cells.Item(19, 17).Value = [double]"206.7479077006513"
$ws.Cells.Item(19, 18).Value = [double]"1240.487446203908"
$ws.Cells.Item(19, 19).Value = [double]"0.005029912539841528"
$ws.Cells.Item(19, 20).Value = [double]"0.003593836661426688"

# Row 20
$ws.Cells.Item(20, 1).Value = "MuSCs"
$ws.Cells.Item(20, 2).Value = "Adam12"
$ws.Cells.Item(20, 3).Value = "Itgb1"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = [double]"2.4758755"
$ws.Cells.Item(20, 8).Value = [double]"4.951751"
$ws.Cells.Item(20, 9).Value = [double]"0.03218061767783864"
$ws.Cells.Item(20, 10).Value = [double]"0.02168637206800571"
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = [double]"91.06846250000001"
$ws.Cells.Item(20, 14).Value = [double]"182.136925"
$ws.Cells.Item(20, 15).Value = [double]"0.1704597085236707"
$ws.Cells.Item(20, 16).Value = [double]"0.1204857969594293"
$ws.Cells.Item(20, 17).Value = [double]"225.4741751264188"
$ws.Cells.Item(20, 18).Value = [double]"901.8967005056751"
$ws.Cells.Item(20, 19).Value = [double]"0.00548549870947606"
$ws.Cells.Item(20, 20).Value = [double]"0.002612899821772376"

# Row 21
$ws.Cells.Item(21, 1).Value = "MuSCs"
$ws.Cells.Item(21, 2).Value = "Adam12"
$ws.Cells.Item(21, 3).Value = "Itgb1"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = [double]"2.4758755"
$ws.Cells.Item(21, 8).Value = [double]"4.951751"
$ws.Cells.Item(21, 9).Value = [double]"0.03218061767783864"
$ws.Cells.Item(21, 10).Value = [double]"0.02168637206800571"
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = [double]"89.83563"
$ws.Cells.Item(21, 14).Value = [double]"269.50689"
$ws.Cells.Item(21, 15).Value = [double]"0.1681521229683693"
$ws.Cells.Item(21, 16).Value = [double]"0.1782820942415013"
$ws.Cells.Item(21, 17).Value = [double]"222.421835344065"
$ws.Cells.Item(21, 18).Value = [double]"1334.53101206439"
$ws.Cells.Item(21, 19).Value = [double]"0.005411239180962001"
$ws.Cells.Item(21, 20).Value = [double]"0.003866291828784456"

# Row 22
$ws.Cells.Item(22, 1).Value = "Resolving-Mac"
$ws.Cells.Item(22, 2).Value = "Adam12"
$ws.Cells.Item(22, 3).Value = "Itgb1"
$ws.Cells.Item(22, 4).Value = "ECs"
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(22, 7).Value = [double]"0.005290666666666667"
$ws.Cells.Item(22, 8).Value = [double]"0.015872"
$ws.Cells.Item(22, 9).Value = [double]"6.876635002886251e-05"
$ws.Cells.Item(22, 10).Value = [double]"6.951199635510484e-05"
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = [double]"121.928739"
$ws.Cells.Item(22, 14).Value = [double]"365.786217"
$ws.Cells.Item(22, 15).Value = [double]"0.2282232151508951"
$ws.Cells.Item(22, 16).Value = [double]"0.2419720431319445"
$ws.Cells.Item(22, 17).Value = [double]"0.645084315136"
$ws.Cells.Item(22, 18).Value = [double]"5.805758836223999"
$ws.Cells.Item(22, 19).Value = [double]"1.569407749777885e-05"
$ws.Cells.Item(22, 20).Value = [double]"1.6819959780225e-05"

# Row 23
$ws.Cells.Item(23, 1).Value = "Resolving-Mac"
$ws.Cells.Item(23, 2).Value = "Adam12"
$ws.Cells.Item(23, 3).Value = "Itgb1"
$ws.Cells.Item(23, 4).Value = "FAPs"
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(23, 7).Value = [double]"0.005290666666666667"
$ws.Cells.Item(23, 8).Value = [double]"0.015872"
$ws.Cells.Item(23, 9).Value = [double]"6.876635002886251e-05"
$ws.Cells.Item(23, 10).Value = [double]"6.951199635510484e-05"
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = [double]"147.91433"
$ws.Cells.Item(23, 14).Value = [double]"443.74299"
$ws.Cells.Item(23, 15).Value = [double]"0.2768624053389947"
$ws.Cells.Item(23, 16).Value = [double]"0.2935413991166814"
$ws.Cells.Item(23, 17).Value = [double]"0.7825654152533332"
$ws.Cells.Item(23, 18).Value = [double]"7.04308873728"
$ws.Cells.Item(23, 19).Value = [double]"1.903881707537412e-05"
$ws.Cells.Item(23, 20).Value = [double]"2.040464866547114e-05"

# Row 24
$ws.Cells.Item(24, 1).Value = "Resolving-Mac"
$ws.Cells.Item(24, 2).Value = "Adam12"
$ws.Cells.Item(24, 3).Value = "Itgb1"
$ws.Cells.Item(24, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(24, 7).Value = [double]"0.005290666666666667"
$ws.Cells.Item(24, 8).Value = [double]"0.015872"
$ws.Cells.Item(24, 9).Value = [double]"6.876635002886251e-05"
$ws.Cells.Item(24, 10).Value = [double]"6.951199635510484e-05"
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = [double]"83.50496933333334"
$ws.Cells.Item(24, 14).Value = [double]"250.514908"
$ws.Cells.Item(24, 15).Value = [double]"0.1563025480180701"
$ws.Cells.Item(24, 16).Value = [double]"0.1657186665504434"
$ws.Cells.Item(24, 17).Value = [double]"0.4417969577528889"
$ws.Cells.Item(24, 18).Value = [double]"3.976172619776"
$ws.Cells.Item(24, 19).Value = [double]"1.07483557274137e-05"
$ws.Cells.Item(24, 20).Value = [double]"1.151943534522726e-05"

# Row 25
$ws.Cells.Item(25, 1).Value = "Resolving-Mac"
$ws.Cells.Item(25, 2).Value = "Adam12"
$ws.Cells.Item(25, 3).Value = "Itgb1"
$ws.Cells.Item(25, 4).Value = "MuSCs"
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(25, 7).Value = [double]"0.005290666666666667"
$ws.Cells.Item(25, 8).Value = [double]"0.015872"
$ws.Cells.Item(25, 9).Value = [double]"6.876635002886251e-05"
$ws.Cells.Item(25, 10).Value = [double]"6.951199635510484e-05"
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = [double]"91.06846250000001"
$ws.Cells.Item(25, 14).Value = [double]"182.136925"
$ws.Cells.Item(25, 15).Value = [double]"0.1704597085236707"
$ws.Cells.Item(25, 16).Value = [double]"0.1204857969594293"
$ws.Cells.Item(25, 17).Value = [double]"0.4818128789333334"
$ws.Cells.Item(25, 18).Value = [double]"2.890877273600001"
$ws.Cells.Item(25, 19).Value = [double]"1.172189198215662e-05"
$ws.Cells.Item(25, 20).Value = [double]"8.375208279085754e-06"

# Row 26
$ws.Cells.Item(26, 1).Value = "Resolving-Mac"
$ws.Cells.Item(26, 2).Value = "Adam12"
$ws.Cells.Item(26, 3).Value = "Itgb1"
$ws.Cells.Item(26, 4).Value = "Resolving-Mac"
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(26, 7).Value = [double]"0.005290666666666667"
$ws.Cells.Item(26, 8).Value = [double]"0.015872"
$ws.Cells.Item(26, 9).Value = [double]"6.876635002886251e-05"
$ws.Cells.Item(26, 10).Value = [double]"6.951199635510484e-05"
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = [double]"89.83563"
$ws.Cells.Item(26, 14).Value = [double]"269.50689"
$ws.Cells.Item(26, 15).Value = [double]"0.1681521229683693"
$ws.Cells.Item(26, 16).Value = [double]"0.1782820942415013"
$ws.Cells.Item(26, 17).Value = [double]"0.4752903731199999"
$ws.Cells.Item(26, 18).Value = [double]"4.27761335808"
$ws.Cells.Item(26, 19).Value = [double]"1.156320774613921e-05"
$ws.Cells.Item(26, 20).Value = [double]"1.23927442850957e-05"
